$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 291.6
$ws.Range("I33").Value = 291.6
$ws.Range("K33").Value = 291.6
$ws.Range("M33").Value = -62.60000000000002

$ws.Range("H100").Value = 1284.3334
$ws.Range("I100").Value = 1216.5385
$ws.Range("K100").Value = 1216.5385
$ws.Range("M100").Value = -675.5385000000001

$ws.Range("H132").Value = 6467.6816
$ws.Range("I132").Value = 6939.45
$ws.Range("K132").Value = 20818.35
$ws.Range("M132").Value = -18288.35

$ws.Range("H138").Value = 1434.1052
$ws.Range("I138").Value = 859.25
$ws.Range("K138").Value = 2577.75
$ws.Range("M138").Value = 2562.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2804.88
$ws.Range("I61").Value = 1947.2
$ws.Range("J61").Value = 4091.4
$ws.Range("K61").Value = 1947.2
$ws.Range("L61").Value = 4091.4
$ws.Range("M61").Value = -1735.2
$ws.Range("N61").Value = -4515.4

$ws.Range("H63").Value = 1498.5
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 1498.5
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 160537.55
$ws.Range("I74").Value = 223592.6
$ws.Range("K74").Value = 223592.6
$ws.Range("M74").Value = -222718.6

$ws.Range("H77").Value = 160537.55
$ws.Range("I77").Value = 223592.6
$ws.Range("K77").Value = 1117963
$ws.Range("M77").Value = -1113595

$ws.Range("H88").Value = 4876.5557
$ws.Range("I88").Value = 2447.5
$ws.Range("J88").Value = 6819.8
$ws.Range("K88").Value = 2447.5
$ws.Range("L88").Value = 6819.8
$ws.Range("M88").Value = -2041.5
$ws.Range("N88").Value = -7631.8

$ws.Range("H91").Value = 4876.5557
$ws.Range("I91").Value = 2447.5
$ws.Range("J91").Value = 6819.8
$ws.Range("K91").Value = 2447.5
$ws.Range("L91").Value = 6819.8
$ws.Range("M91").Value = -1043.5
$ws.Range("N91").Value = -9627.799999999999

$ws.Range("H97").Value = 584.16
$ws.Range("I97").Value = 556.4167
$ws.Range("K97").Value = 556.4167
$ws.Range("M97").Value = -60.41669999999999

$ws.Range("H132").Value = 4115.8335
$ws.Range("I132").Value = 4642.8184
$ws.Range("K132").Value = 13928.4552
$ws.Range("M132").Value = -11398.4552

$ws.Range("H136").Value = 2804.88
$ws.Range("I136").Value = 1947.2
$ws.Range("J136").Value = 4091.4
$ws.Range("K136").Value = 5841.6
$ws.Range("L136").Value = 12274.2
$ws.Range("M136").Value = -3291.6
$ws.Range("N136").Value = -17374.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 6500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 6500
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -8746

$ws.Range("H89").Value = 6500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 6500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 32500
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -43732

$ws.Range("H99").Value = 3810.625
$ws.Range("I99").Value = 3996.6667
$ws.Range("K99").Value = 3996.6667
$ws.Range("M99").Value = -2498.6667

$ws.Range("H107").Value = 2263808.8
$ws.Range("I107").Value = 2653823
$ws.Range("K107").Value = 2653823
$ws.Range("M107").Value = -2651903

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2027.7273
$ws.Range("I16").Value = 2034.3334
$ws.Range("J16").Value = 2019.8
$ws.Range("K16").Value = 2034.3334
$ws.Range("L16").Value = 2019.8
$ws.Range("M16").Value = -1747.3334
$ws.Range("N16").Value = -2593.8

$ws.Range("H22").Value = 341.83334
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -1000

$ws.Range("H58").Value = 3166.6667
$ws.Range("I58").Value = 2500
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 2500
$ws.Range("L58").Value = 3500
$ws.Range("M58").Value = -2297
$ws.Range("N58").Value = -3906

$ws.Range("H113").Value = 2027.7273
$ws.Range("I113").Value = 2034.3334
$ws.Range("J113").Value = 2019.8
$ws.Range("K113").Value = 2034.3334
$ws.Range("L113").Value = 2019.8
$ws.Range("M113").Value = 135.6666
$ws.Range("N113").Value = -6359.8

$ws.Range("H132").Value = 4027.8948
$ws.Range("I132").Value = 3915.8333
$ws.Range("J132").Value = 4220
$ws.Range("K132").Value = 11747.4999
$ws.Range("L132").Value = 12660
$ws.Range("M132").Value = -9217.499899999999
$ws.Range("N132").Value = -17720

$ws.Range("H136").Value = 3166.6667
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5005171.5
$ws.Range("J68").Value = 11118632
$ws.Range("L68").Value = 33355896
$ws.Range("N68").Value = -33357518

$ws.Range("H71").Value = 5005171.5
$ws.Range("J71").Value = 11118632
$ws.Range("L71").Value = 100067688
$ws.Range("N71").Value = -100075800

$ws.Range("H74").Value = 25711.572
$ws.Range("J74").Value = 25993.6
$ws.Range("L74").Value = 77980.79999999999
$ws.Range("N74").Value = -80102.79999999999

$ws.Range("H77").Value = 25711.572
$ws.Range("J77").Value = 25993.6
$ws.Range("L77").Value = 233942.4
$ws.Range("N77").Value = -244550.4

$ws.Range("H132").Value = 2037.6
$ws.Range("J132").Value = 2226.25
$ws.Range("L132").Value = 20036.25
$ws.Range("N132").Value = -25096.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2371.1292
$ws.Range("I132").Value = 1973.5385
$ws.Range("J132").Value = 2658.2778
$ws.Range("K132").Value = 5920.6155
$ws.Range("L132").Value = 7974.8334
$ws.Range("M132").Value = -3390.6155
$ws.Range("N132").Value = -13034.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6876.2
$ws.Range("J7").Value = 9250
$ws.Range("L7").Value = 9250
$ws.Range("N7").Value = -9474

$ws.Range("H22").Value = 1653.5834
$ws.Range("I22").Value = 1707.4
$ws.Range("J22").Value = 1384.5
$ws.Range("K22").Value = 1707.4
$ws.Range("L22").Value = 1384.5
$ws.Range("M22").Value = -1412.4
$ws.Range("N22").Value = -1974.5

$ws.Range("H27").Value = 1653.5834
$ws.Range("I27").Value = 1707.4
$ws.Range("J27").Value = 1384.5
$ws.Range("K27").Value = 1707.4
$ws.Range("L27").Value = 1384.5
$ws.Range("M27").Value = -1600.4
$ws.Range("N27").Value = -1598.5

$ws.Range("H126").Value = 6876.2
$ws.Range("J126").Value = 9250
$ws.Range("L126").Value = 27750
$ws.Range("N126").Value = -32690

$ws.Range("H132").Value = 4778.394
$ws.Range("I132").Value = 4216
$ws.Range("J132").Value = 5192.7896
$ws.Range("K132").Value = 12648
$ws.Range("L132").Value = 15578.3688
$ws.Range("M132").Value = -10118
$ws.Range("N132").Value = -20638.3688

$ws.Range("H136").Value = 6131.6665
$ws.Range("J136").Value = 8000
$ws.Range("L136").Value = 24000
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8832.666999999999
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 8832.666999999999
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 8832.666999999999
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -10080.667

$ws.Range("H65").Value = 8832.666999999999
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 8832.666999999999
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 44163.335
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -50403.335

$ws.Range("H81").Value = 5564.6665
$ws.Range("I81").Value = 6683.7144
$ws.Range("J81").Value = 3998
$ws.Range("K81").Value = 13367.4288
$ws.Range("L81").Value = 7996
$ws.Range("M81").Value = -12306.4288
$ws.Range("N81").Value = -10118

$ws.Range("H84").Value = 5564.6665
$ws.Range("I84").Value = 6683.7144
$ws.Range("J84").Value = 3998
$ws.Range("K84").Value = 66837.144
$ws.Range("L84").Value = 39980
$ws.Range("M84").Value = -61533.144
$ws.Range("N84").Value = -50588

$ws.Range("H122").Value = 31252888
$ws.Range("I122").Value = 3184.8333
$ws.Range("K122").Value = 9554.499899999999
$ws.Range("M122").Value = -7104.499899999999

$ws.Range("H132").Value = 1297.5834
$ws.Range("I132").Value = 1201.8928
$ws.Range("K132").Value = 3605.6784
$ws.Range("M132").Value = -1075.6784

$ws.Range("H136").Value = 9248.9
$ws.Range("I136").Value = 4720.4443
$ws.Range("K136").Value = 14161.3329
$ws.Range("M136").Value = -11611.3329
